# ---------------------------------------------------------------------------
# Target change (from the supplied OOXML diff):
#
#   word/numbering.xml -- four <w:abstractNum> definitions (abstractNumId
#   990, 991, 99411, 99414) each get a new, unrelated-looking random value
#   for their <w:nsid w:val="..."/> child element. Nothing else in those
#   abstractNum blocks (multiLevelType, lvl definitions, numFmt, lvlText,
#   indentation, ...) changes, and no <w:num>/<w:abstractNumId> mapping,
#   paragraph, or numPr reference anywhere in document.xml changes either.
#
# Investigation performed before writing this script (see transcript):
#   * <w:nsid> is Word's internal "list signature" GUID fragment. It is
#     generated by Word when a list definition is minted and is used
#     purely to detect/merge duplicate list definitions when documents
#     are combined. Exactly like in real Microsoft Word, it is NOT part
#     of the Word object model exposed to VBA/COM automation: there is
#     no ListTemplate.Nsid, List.Nsid, ListFormat.Nsid, etc. anywhere in
#     the object model (checked the full Document.list_commands() surface
#     - 4433 Class.Member entries - plus every List*/Num* class; none of
#     them expose it).
#   * Document.WordOpenXML and Range.XML are get-only in this host (the
#     runtime raises "... is a read-only property" if you try to assign
#     them), so the usual "round-trip the whole package as a string"
#     trick is unavailable; the only sanctioned way to change content is
#     through Range.Text / InsertXML on the exact body Range whose
#     *content* should change - and <w:abstractNum> lives in the separate
#     numbering part, not in any body Range.
#   * Find.Execute only searches visible document text, never attribute
#     values inside part XML, so it cannot reach <w:nsid>.
#   * Methods that do mutate numbering state (ListFormat.ApplyListTemplate,
#     ListFormat.StartNewList, ListFormat.RemoveNumbers) were exercised
#     directly: StartNewList/RemoveNumbers never touch numbering.xml at
#     all, and ApplyListTemplate only ever *mints a brand-new*
#     abstractNumId (appending, e.g. 99415, 99416, ...) - real Word never
#     rewrites an existing abstractNum's <w:nsid> in place, and neither
#     does this host. Forcing that path would also rewrite every
#     paragraph's <w:numPr>/numId throughout the body (152 list
#     paragraphs here), which is a far larger and incorrect divergence
#     from the target than making no change at all.
#   * ListLevel.* / ListTemplate.* property setters (NumberFormat,
#     Alignment, TrailingCharacter, StartAt, Name, OutlineNumbered, ...)
#     were probed too: assigning them does not bump the document's
#     mutation generation counter and never changes a single byte of the
#     saved numbering.xml, confirming they are read-mostly/cosmetic in
#     this host and provide no write path into <w:abstractNum>.
#
# Conclusion: nothing reachable through Word.Application / Document COM
# automation (real or emulated) can rewrite <w:nsid>, which is a random,
# non-content-bearing identifier Word itself assigns internally. The
# commit message ("changed broken link paths") does not correspond to any
# visible edit in this document - the four nsid values are incidental
# churn from whatever external process regenerated this file, not a
# change an end user/macro could have made in Word. Reproducing it here
# would require writing raw package bytes, which this automation surface
# intentionally does not allow (by design, matching real Word).
#
# So: touch nothing. Re-saving the document unmodified is the closest,
# and only non-destructive, reachable approximation of the target state -
# every other avenue tested actively makes the document diverge further
# from the target (minted list definitions, rewritten numId references,
# etc.) without ever managing to touch the actual <w:nsid> values.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# No-op: confirm we still have the document in hand without mutating it.
$null = $d.Name
